$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency "Price" values in column D look numeric (e.g. "578.32" or
# "67.382.05") but must stay verbatim TEXT, matching the sheet author intent
# (the workbook stores them as inline strings, not numbers). Writing the raw
# string through .Value lets Excel "smart" auto-detect them as numbers, which
# both changes the cell type and mangles values like "578.32" into a binary
# float (578.32000000000005) on save. Forcing a Text number format before the
# write - then restoring the default "Normal" style right after - keeps the
# value as exact text without leaving any stray number formatting behind.
function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "67.382.05"
$ws.Range("E2").Value = "  -1.23%  "
Set-TextValue "D3" "3.219.23"
$ws.Range("E3").Value = "  -1.61%  "
Set-TextValue "D5" "578.32"
$ws.Range("E5").Value = "  -1.51%  "
Set-TextValue "D6" "182.14"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.19%  "
Set-TextValue "D9" "3.216.24"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -3.47%  "
Set-TextValue "D11" "6.58"
$ws.Range("E11").Value = "  -1.96%  "
Set-TextValue "D12" "0.411"
$ws.Range("E12").Value = "  -1.67%  "
Set-TextValue "D13" "3.779.98"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("E14").Value = "  -0.17%  "
Set-TextValue "D15" "27.64"
$ws.Range("E15").Value = "  -3.66%  "
Set-TextValue "D16" "67.443.17"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("E17").Value = "  -2.70%  "
Set-TextValue "D18" "3.238.74"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -2.36%  "
Set-TextValue "D20" "13.38"
$ws.Range("E20").Value = "  -1.89%  "
Set-TextValue "D21" "394.84"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -0.16%  "
Set-TextValue "D24" "70.71"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("E27").Value = "  +0.17%  "
Set-TextValue "D28" "9.52"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -2.39%  "
Set-TextValue "D31" "5.55"
$ws.Range("E31").Value = "  -3.82%  "
Set-TextValue "D32" "22.59"
$ws.Range("E32").Value = "  -1.65%  "
Set-TextValue "D33" "6.98"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("E34").Value = "  -0.04%  "
Set-TextValue "D35" "1.24"
$ws.Range("E35").Value = "  -3.17%  "
Set-TextValue "D36" "160.94"
$ws.Range("E36").Value = "  -1.09%  "
Set-TextValue "D37" "1.46"
$ws.Range("E37").Value = "  -5.58%  "
Set-TextValue "D38" "1.87"
$ws.Range("E38").Value = "  -0.16%  "
Set-TextValue "D39" "26.23"
$ws.Range("E39").Value = "  -1.91%  "
Set-TextValue "D40" "0.802"
$ws.Range("E40").Value = "  -4.22%  "
Set-TextValue "D41" "4.55"
$ws.Range("E41").Value = "  -1.36%  "
Set-TextValue "D42" "6.48"
$ws.Range("E42").Value = "  -4.49%  "
Set-TextValue "D43" "2.46"
$ws.Range("E43").Value = "  -6.27%  "
$ws.Range("E44").Value = "  -1.33%  "
Set-TextValue "D45" "40.51"
$ws.Range("E45").Value = "  -2.44%  "
Set-TextValue "D46" "2.603.31"
Set-TextValue "D47" "24.44"
$ws.Range("E47").Value = "  -4.16%  "
Set-TextValue "D48" "332.47"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -1.97%  "
